$wb = $excel.ActiveWorkbook

# Work on sheet "M"
$ws = $wb.Worksheets.Item("M")

# Update the two shared-string header values on row 1
$ws.Range("B1").Value = "UMAR-SURS--MZ002--HR--M"
$ws.Range("C1").Value = "UMAR-SURS--MZ002--SI--M"

# Remove column D (header + all data) -- data now ends at column C
$ws.Columns.Item(4).Delete()

# Set explicit width on column A (value of 12.1 characters rounds/stores as 13)
$ws.Columns.Item(1).ColumnWidth = 12.1

# Move the active selection
$ws.Range("F10").Select()

# Remove the second worksheet "A"
$wsA = $wb.Worksheets.Item("A")
$wsA.Delete()
